$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "datos actualizados" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 21:35"

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4,2).Value = 1609921
$ws.Cells.Item(4,3).Value = 17198
$ws.Cells.Item(4,4).Value = 373638
$ws.Cells.Item(4,5).Value = 1140444
$ws.Cells.Item(4,7).Value = 903
$ws.Cells.Item(4,8).Value = 95839

# --- Row 6: Brasil ---
$ws.Cells.Item(6,2).Value = 296113
$ws.Cells.Item(6,3).Value = 2756
$ws.Cells.Item(6,5).Value = 160274
$ws.Cells.Item(6,7).Value = 262
$ws.Cells.Item(6,8).Value = 19156

# --- Row 11: Alemania ---
$ws.Cells.Item(11,2).Value = 178886
$ws.Cells.Item(11,3).Value = 355
$ws.Cells.Item(11,5).Value = 12605
$ws.Cells.Item(11,7).Value = 11
$ws.Cells.Item(11,8).Value = 8281

# --- Row 113: Costa Rica ---
$ws.Cells.Item(113,2).Value = 903
$ws.Cells.Item(113,3).Value = 6
$ws.Cells.Item(113,4).Value = 592
$ws.Cells.Item(113,5).Value = 301

# --- Row 141: Togo ---
$ws.Cells.Item(141,2).Value = 354
$ws.Cells.Item(141,3).Value = 14
$ws.Cells.Item(141,4).Value = 118
$ws.Cells.Item(141,5).Value = 224

# --- Row 147: Nicaragua ---
$ws.Cells.Item(147,2).Value = 279
$ws.Cells.Item(147,3).Value = 25
$ws.Cells.Item(147,5).Value = 63

# --- Rows 156-160: Mauritania re-enters the top-220 table (pushing Mozambique,
#     Guadalupe, Gibraltar and Uganda down one row each) ---
$ws.Cells.Item(156,1).Value = "Mauritania"
$ws.Cells.Item(156,2).Value = 173
$ws.Cells.Item(156,3).Value = 32
$ws.Cells.Item(156,4).Value = 7
$ws.Cells.Item(156,5).Value = 162
$ws.Cells.Item(156,8).Value = 4

$ws.Cells.Item(157,1).Value = "Mozambique"
$ws.Cells.Item(157,2).Value = 162
$ws.Cells.Item(157,3).Value = 6
$ws.Cells.Item(157,4).Value = 48
$ws.Cells.Item(157,5).Value = 114
$ws.Cells.Item(157,8).Value = 0

$ws.Cells.Item(158,1).Value = "Guadalupe"
$ws.Cells.Item(158,2).Value = 155
$ws.Cells.Item(158,3).Value = 0
$ws.Cells.Item(158,4).Value = 109
$ws.Cells.Item(158,5).Value = 33
$ws.Cells.Item(158,8).Value = 13

$ws.Cells.Item(159,1).Value = "Gibraltar"
$ws.Cells.Item(159,2).Value = 151
$ws.Cells.Item(159,3).Value = 2
$ws.Cells.Item(159,4).Value = 146
$ws.Cells.Item(159,5).Value = 5

$ws.Cells.Item(160,1).Value = "Uganda"
$ws.Cells.Item(160,2).Value = 145
$ws.Cells.Item(160,4).Value = 65
$ws.Cells.Item(160,5).Value = 80
$ws.Cells.Item(160,8).Value = 0

# --- Row 175: Malaui ---
$ws.Cells.Item(175,2).Value = 72
$ws.Cells.Item(175,3).Value = 1
$ws.Cells.Item(175,5).Value = 42
